# Update the MAIN DASHBOARD sheet with revised energy/exposure figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=HOUR, B=ACTUAL_ENERGY, C=CONTESTABLE_ENERGY,
#          D=TOTAL_BCQ_NOMINATION (unchanged), E=WESM_EXPOSURE
# Data rows run from row 2 (HOUR 1) through row 25 (HOUR 24).

$data = @(
    @{ Row = 2;  B = 73347.2945;          C = 5623.7055;            E = 2723.589000000007 }
    @{ Row = 3;  B = 70244.77099999999;   C = 5458.229;             E = 19786.54199999999 }
    @{ Row = 4;  B = 67129.0255;          C = 5507.9745;            E = 16621.05100000001 }
    @{ Row = 5;  B = 64595.502;           C = 5475.498;             E = 14120.004 }
    @{ Row = 6;  B = 65374.48;            C = 5469.52;              E = 14904.96000000001 }
    @{ Row = 7;  B = 67982.29700000001;   C = 5575.703;             E = 17406.594 }
    @{ Row = 8;  B = 67367.1265;          C = 6141.8735;            E = 38725.253 }
    @{ Row = 9;  B = 79207.3645;          C = 7163.6355;            E = 27043.72899999999 }
    @{ Row = 10; B = 95174.5395;          C = 8647.460500000001;    E = 21527.079 }
    @{ Row = 11; B = 81952.2095;          C = 12435.7905;           E = 4516.418999999994 }
    @{ Row = 12; B = 81626.978;           C = 15471.022;            E = 1155.956000000006 }
    @{ Row = 13; B = 80925.1525;          C = 16384.8475;           E = -459.695000000007 }
    @{ Row = 14; B = 80952.772;           C = 16062.228;            E = -109.4560000000056 }
    @{ Row = 15; B = 107586.796;          C = 16185.204;            E = 26401.592 }
    @{ Row = 16; B = 106447.456;          C = 16525.544;            E = 24921.91200000001 }
    @{ Row = 17; B = 65630.38250000001;   C = 16462.6175;           E = -15832.23499999999 }
    @{ Row = 18; B = 40052.2725;          C = 16761.7275;           E = -41709.455 }
    @{ Row = 19; C = 16532.1625;          E = 7051.228499999997 }
    @{ Row = 20; C = 15890.3885;          E = 8134.193499999994 }
    @{ Row = 21; C = 13683.775;           E = 8380.1875 }
    @{ Row = 22; C = 11966.283;           E = 8944.044999999998 }
    @{ Row = 23; C = 9398.494000000001;   E = 10285.514 }
    @{ Row = 24; C = 6574.9635;           E = 10044.24250000001 }
    @{ Row = 25; C = 5831.934499999999;   E = -18483.906 }
)

foreach ($item in $data) {
    $r = $item.Row
    if ($item.ContainsKey('B')) {
        $ws.Cells.Item($r, 2).Value = $item.B
    }
    if ($item.ContainsKey('C')) {
        $ws.Cells.Item($r, 3).Value = $item.C
    }
    if ($item.ContainsKey('E')) {
        $ws.Cells.Item($r, 5).Value = $item.E
    }
}
